# Update cryptocurrency price (D) and 1h volume change (E) columns
# for rows 2-51 on the active worksheet, per the latest GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.270.17'
$ws.Range("E2").Value = '  +1.76%  '
$ws.Range("D3").Value = '2.629.31'
$ws.Range("E3").Value = '  +0.05%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '601.89'
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("D6").Value = '151.15'
$ws.Range("E6").Value = '  +3.27%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = "'0.590"
$ws.Range("E8").Value = '  +0.76%  '
$ws.Range("E9").Value = '  +2.29%  '
$ws.Range("E10").Value = '  +2.99%  '
$ws.Range("E11").Value = '  +6.84%  '
$ws.Range("E12").Value = '  -0.81%  '
$ws.Range("D13").Value = '27.72'
$ws.Range("E13").Value = '  +2.10%  '
$ws.Range("D14").Value = '3.099.06'
$ws.Range("E14").Value = '  +0.03%  '
$ws.Range("D15").Value = '64.103.92'
$ws.Range("E15").Value = '  +1.71%  '
$ws.Range("E16").Value = '  +3.94%  '
$ws.Range("D17").Value = '2.625.99'
$ws.Range("E17").Value = '  +2.64%  '
$ws.Range("D18").Value = '12.18'
$ws.Range("E18").Value = '  +7.94%  '
$ws.Range("D19").Value = '4.67'
$ws.Range("E19").Value = '  +3.91%  '
$ws.Range("D20").Value = '353.01'
$ws.Range("E20").Value = '  +3.88%  '
$ws.Range("D21").Value = '6.98'
$ws.Range("E21").Value = '  +1.29%  '
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").Value = '5.72'
$ws.Range("E23").Value = '  +2.99%  '
$ws.Range("D24").Value = "'66.70"
$ws.Range("E24").Value = '  +0.48%  '
$ws.Range("E25").Value = '  +15.77%  '
$ws.Range("E26").Value = '  +5.12%  '
$ws.Range("D27").Value = '9.29'
$ws.Range("E27").Value = '  +7.43%  '
$ws.Range("E28").Value = '  +2.09%  '
$ws.Range("D29").Value = "'8.10"
$ws.Range("E29").Value = '  +3.37%  '
$ws.Range("D30").Value = "'542.30"
$ws.Range("E30").Value = '  +0.40%  '
$ws.Range("E31").Value = '  +0.04%  '
$ws.Range("D32").Value = '2.06'
$ws.Range("E32").Value = '  +2.10%  '
$ws.Range("D33").Value = '0.0₃0857'
$ws.Range("E33").Value = '  +7.02%  '
$ws.Range("D34").Value = '1.75'
$ws.Range("E34").Value = '  +0.55%  '
$ws.Range("D35").Value = '5.27'
$ws.Range("E35").Value = '  -1.08%  '
$ws.Range("D36").Value = "'167.30"
$ws.Range("E36").Value = '  -0.05%  '
$ws.Range("E37").Value = '  +7.07%  '
$ws.Range("E38").Value = '  +2.00%  '
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("D40").Value = '19.58'
$ws.Range("E40").Value = '  +3.16%  '
$ws.Range("D41").Value = '171.28'
$ws.Range("E41").Value = '  +1.76%  '
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("D43").Value = '40.04'
$ws.Range("E43").Value = '  +0.83%  '
$ws.Range("E44").Value = '  +5.39%  '
$ws.Range("D45").Value = "'0.0590"
$ws.Range("E45").Value = '  +4.46%  '
$ws.Range("D46").Value = '21.65'
$ws.Range("E46").Value = '  -2.88%  '
$ws.Range("E47").Value = '  +1.22%  '
$ws.Range("E48").Value = '  +14.79%  '
$ws.Range("D49").Value = '0.0247'
$ws.Range("E49").Value = '  +2.23%  '
$ws.Range("D50").Value = '0.0967'
$ws.Range("E50").Value = '  +1.09%  '
$ws.Range("E51").Value = '  +4.10%  '
